# Generate Report for Handoff
#
# A new handoff xliff was generated for the
# "6113a586-8ee9-4b80-892d-43107333790f.md" file, so its
# "Latest Handoff Datetime" column (H) is refreshed on both the
# zh-cn and de-de localization-status worksheets.

$wb = $excel.ActiveWorkbook

# zh-cn sheet (row 4 corresponds to the 6113a586 file)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-12-15 03:40:37"

# de-de sheet (row 4 corresponds to the 6113a586 file)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-12-15 03:40:49"
